$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8 (G=4565)
$ws.Range("H8").Value = 53.2
$ws.Range("I8").Value = 65.375
$ws.Range("J8").Value = 4.5
$ws.Range("K8").Value = 196.125
$ws.Range("L8").Value = 13.5
$ws.Range("M8").Value = -57.125
$ws.Range("N8").Value = -291.5

# Row 17 (G=38956)
$ws.Range("H17").Value = 1268
$ws.Range("I17").Value = 900
$ws.Range("J17").Value = 1360
$ws.Range("K17").Value = 2700
$ws.Range("L17").Value = 4080
$ws.Range("M17").Value = -2532
$ws.Range("N17").Value = -4416

# Row 43 (G=5472)
$ws.Range("H43").Value = 77142.42999999999
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 77142.42999999999
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 77142.42999999999
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -77280.42999999999

# Row 51 (G=5486)
$ws.Range("H51").Value = 5000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 5000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 5000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -5968

# Row 106 (G=19903)
$ws.Range("H106").Value = 5000
$ws.Range("I106").Value = 5000
$ws.Range("K106").Value = 5000
$ws.Range("M106").Value = -4369

# Row 107 (G=27766)
$ws.Range("H107").Value = 2260.75
$ws.Range("I107").Value = 2260.75
$ws.Range("K107").Value = 2260.75
$ws.Range("M107").Value = -340.75

# Row 129 (G=36115)
$ws.Range("H129").Value = 964.375
$ws.Range("I129").Value = 964.375
$ws.Range("K129").Value = 2893.125
$ws.Range("M129").Value = 2106.875

# Row 138 (G=44169)
$ws.Range("H138").Value = 2888
$ws.Range("I138").Value = 2418.6667
$ws.Range("K138").Value = 7256.000100000001
$ws.Range("M138").Value = -2116.000100000001

$ws = $wb.Worksheets.Item("ARM")
# Row 74 (G=44000)
$ws.Range("H74").Value = 2399.75
$ws.Range("I74").Value = 2399.75
$ws.Range("K74").Value = 2399.75
$ws.Range("M74").Value = -1525.75

# Row 77 (G=44000)
$ws.Range("H77").Value = 2399.75
$ws.Range("I77").Value = 2399.75
$ws.Range("K77").Value = 11998.75
$ws.Range("M77").Value = -7630.75

# Row 88 (G=12530)
$ws.Range("H88").Value = 3281.238
$ws.Range("I88").Value = 1542.6
$ws.Range("J88").Value = 3824.5625
$ws.Range("K88").Value = 1542.6
$ws.Range("L88").Value = 3824.5625
$ws.Range("M88").Value = -1136.6
$ws.Range("N88").Value = -4636.5625

# Row 91 (G=12530)
$ws.Range("H91").Value = 3281.238
$ws.Range("I91").Value = 1542.6
$ws.Range("J91").Value = 3824.5625
$ws.Range("K91").Value = 1542.6
$ws.Range("L91").Value = 3824.5625
$ws.Range("M91").Value = -138.5999999999999
$ws.Range("N91").Value = -6632.5625

$ws = $wb.Worksheets.Item("BSM")
# Row 99 (G=19943)
$ws.Range("H99").Value = 3095.0715
$ws.Range("I99").Value = 2073.2
$ws.Range("K99").Value = 2073.2
$ws.Range("M99").Value = -575.1999999999998

# Row 105 (G=19947)
$ws.Range("H105").Value = 3581.4666
$ws.Range("I105").Value = 3185.25
$ws.Range("K105").Value = 3185.25
$ws.Range("M105").Value = -1438.25

# Row 132 (G=41855)
$ws.Range("H132").Value = 49999
$ws.Range("J132").Value = 49999
$ws.Range("L132").Value = 49999
$ws.Range("N132").Value = -60119

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (G=5367)
$ws.Range("H22").Value = 612.25
$ws.Range("I22").Value = 612.25
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 612.25
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -262.25
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (G=4650)
$ws.Range("H4").Value = 725.25
$ws.Range("I4").Value = 670.75
$ws.Range("K4").Value = 2012.25
$ws.Range("M4").Value = -1900.25

# Row 56 (G=10146)
$ws.Range("H56").Value = 9998.333000000001
$ws.Range("I56").Value = 9998.333000000001
$ws.Range("K56").Value = 9998.333000000001
$ws.Range("M56").Value = -9468.333000000001

# Row 113 (G=27843)
$ws.Range("H113").Value = 950.75
$ws.Range("J113").Value = 950
$ws.Range("L113").Value = 2850
$ws.Range("N113").Value = -7190

# Row 131 (G=36060)
$ws.Range("H131").Value = 1112.909
$ws.Range("I131").Value = 995.6667
$ws.Range("J131").Value = 1131.421
$ws.Range("K131").Value = 2987.0001
$ws.Range("L131").Value = 3394.263
$ws.Range("M131").Value = 2052.9999
$ws.Range("N131").Value = -13474.263

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (G=14146)
$ws.Range("H70").Value = 10000
$ws.Range("I70").Value = 12000
$ws.Range("J70").Value = 8000
$ws.Range("K70").Value = 12000
$ws.Range("L70").Value = 8000
$ws.Range("M70").Value = -11730
$ws.Range("N70").Value = -8540

# Row 73 (G=14146)
$ws.Range("H73").Value = 10000
$ws.Range("I73").Value = 12000
$ws.Range("J73").Value = 8000
$ws.Range("K73").Value = 12000
$ws.Range("L73").Value = 8000
$ws.Range("M73").Value = -11064
$ws.Range("N73").Value = -9872

# Row 133 (G=41854)
$ws.Range("H133").Value = 67186.336
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 67186.336
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 67186.336
$ws.Range("M133").ClearContents()
$ws.Range("N133").Value = -77306.336

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (G=5277)
$ws.Range("H22").Value = 2036.1111
$ws.Range("I22").Value = 1339.8
$ws.Range("J22").Value = 2906.5
$ws.Range("K22").Value = 1339.8
$ws.Range("L22").Value = 2906.5
$ws.Range("M22").Value = -1044.8
$ws.Range("N22").Value = -3496.5

# Row 27 (G=5277)
$ws.Range("H27").Value = 2036.1111
$ws.Range("I27").Value = 1339.8
$ws.Range("J27").Value = 2906.5
$ws.Range("K27").Value = 1339.8
$ws.Range("L27").Value = 2906.5
$ws.Range("M27").Value = -1232.8
$ws.Range("N27").Value = -3120.5

# Row 68 (G=12563)
$ws.Range("H68").Value = 2190.75
$ws.Range("I68").Value = 2190.75
$ws.Range("K68").Value = 2190.75
$ws.Range("M68").Value = -1441.75

# Row 71 (G=12563)
$ws.Range("H71").Value = 2190.75
$ws.Range("I71").Value = 2190.75
$ws.Range("K71").Value = 10953.75
$ws.Range("M71").Value = -7209.75

$ws = $wb.Worksheets.Item("WVR")
# Row 62 (G=12589)
$ws.Range("H62").Value = 2321.1428
$ws.Range("I62").Value = 1609.6
$ws.Range("J62").Value = 4100
$ws.Range("K62").Value = 1609.6
$ws.Range("L62").Value = 4100
$ws.Range("M62").Value = -985.5999999999999
$ws.Range("N62").Value = -5348

# Row 65 (G=12589)
$ws.Range("H65").Value = 2321.1428
$ws.Range("I65").Value = 1609.6
$ws.Range("J65").Value = 4100
$ws.Range("K65").Value = 8048
$ws.Range("L65").Value = 20500
$ws.Range("M65").Value = -4928
$ws.Range("N65").Value = -26740

# Row 132 (G=44029)
$ws.Range("H132").Value = 2460.6365
$ws.Range("I132").Value = 2406.2
$ws.Range("K132").Value = 7218.599999999999
$ws.Range("M132").Value = -4688.599999999999
